$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the list of cells that receive new text values.
$cellRefs = @(
    "D2", "E2", "D3", "E3", "E4", "D5", "E5", "D6", "E6", "E7", "D8", "E8", "D9", "E9", "E10", "D11",
    "E11", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "D18", "E18",
    "E19", "D20", "E21", "E22", "D23", "E23", "D24", "E24", "E25", "E26", "D27", "E27", "D28", "E28",
    "E29", "D30", "E30", "E31", "D32", "E32", "E33", "E34", "E35", "D36", "E36", "E37", "D38", "E38",
    "D39", "E39", "D40", "E40", "E41", "E42", "D43", "E43", "E44", "D45", "E45", "B46", "C46", "D46",
    "E46", "B47", "C47", "D47", "E47", "E48", "D49", "E49", "D50", "E50", "D51", "E51"
)

# Force text format on the cells we are about to rewrite so that Excel
# does not reinterpret numeric-looking strings (e.g. "550.81") as numbers.
foreach ($ref in $cellRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "62.025.21"
$ws.Range("E2").Value = "  -2.39%  "
$ws.Range("D3").Value = "2.576.94"
$ws.Range("E3").Value = "  -4.31%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "550.81"
$ws.Range("E5").Value = "  -1.17%  "
$ws.Range("D6").Value = "155.38"
$ws.Range("E6").Value = "  -2.34%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "0.596"
$ws.Range("E8").Value = "  +1.99%  "
$ws.Range("D9").Value = "0.105"
$ws.Range("E9").Value = "  -1.14%  "
$ws.Range("E10").Value = "  -1.63%  "
$ws.Range("D11").Value = "5.58"
$ws.Range("E11").Value = "  +3.72%  "
$ws.Range("E12").Value = "  -0.89%  "
$ws.Range("D13").Value = "3.033.91"
$ws.Range("E13").Value = "  -4.31%  "
$ws.Range("D14").Value = "25.75"
$ws.Range("E14").Value = "  -2.89%  "
$ws.Range("D15").Value = "61.899.08"
$ws.Range("E15").Value = "  -2.36%  "
$ws.Range("D16").Value = "0.0000145"
$ws.Range("E16").Value = "  +0.05%  "
$ws.Range("D17").Value = "2.581.12"
$ws.Range("E17").Value = "  -4.27%  "
$ws.Range("D18").Value = "11.62"
$ws.Range("E18").Value = "  -4.03%  "
$ws.Range("E19").Value = "  -0.43%  "
$ws.Range("D20").Value = "338.62"
$ws.Range("E21").Value = "  -4.63%  "
$ws.Range("E22").Value = "  +0.39%  "
$ws.Range("D23").Value = "0.494"
$ws.Range("E23").Value = "  -2.76%  "
$ws.Range("D24").Value = "63.48"
$ws.Range("E24").Value = "  -0.78%  "
$ws.Range("E25").Value = "  -0.66%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").Value = "8.18"
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("D28").Value = "7.36"
$ws.Range("E28").Value = "  +4.48%  "
$ws.Range("E29").Value = "  -2.30%  "
$ws.Range("D30").Value = "1.35"
$ws.Range("E30").Value = "  +1.24%  "
$ws.Range("E31").Value = "  -2.67%  "
$ws.Range("D32").Value = "162.74"
$ws.Range("E32").Value = "  -1.88%  "
$ws.Range("E33").Value = "  +1.37%  "
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("E35").Value = "  +0.55%  "
$ws.Range("D36").Value = "19.25"
$ws.Range("E36").Value = "  -1.80%  "
$ws.Range("E37").Value = "  +0.34%  "
$ws.Range("D38").Value = "331.08"
$ws.Range("E38").Value = "  -3.80%  "
$ws.Range("D39").Value = "6.03"
$ws.Range("E39").Value = "  -1.59%  "
$ws.Range("D40").Value = "0.916"
$ws.Range("E40").Value = "  -3.73%  "
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("E42").Value = "  -1.76%  "
$ws.Range("D43").Value = "20.99"
$ws.Range("E43").Value = "  +0.21%  "
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("D45").Value = "0.608"
$ws.Range("E45").Value = "  -2.58%  "
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").Value = "0.0550"
$ws.Range("E46").Value = "  -2.77%  "
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "2.117.82"
$ws.Range("E47").Value = "  +0.47%  "
$ws.Range("E48").Value = "  -0.93%  "
$ws.Range("D49").Value = "19.60"
$ws.Range("E49").Value = "  -3.98%  "
$ws.Range("D50").Value = "0.0967"
$ws.Range("E50").Value = "  -0.72%  "
$ws.Range("D51").Value = "0.0240"
$ws.Range("E51").Value = "  -0.97%  "

# Restore the default (Normal) style so the cells keep their original
# unstyled appearance, matching the source workbook formatting.
foreach ($ref in $cellRefs) {
    $ws.Range($ref).Style = "Normal"
}
